# Tesla numbers update:
#  - rename several header labels (car* -> automotive*/delivered*/produced*)
#  - insert a new "producedCars" column (H+I) before the old "carNumCars"
#    column, which itself is renamed to "deliveredCars" and keeps its
#    original F+G formula.
#  - all columns from the old J (carNumCars) onward shift one column right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at J; everything from J onward (old carNumCars,
# energyRevenue, energyCostOfRevenue, energyStorage) shifts right to
# K, L, M, N respectively. Excel copies the format of the inserted
# column from the cells that get pushed right, matching the diff's
# style attributes.
$ws.Columns.Item(10).Insert()

# Re-label the headers that changed wording (row 1).
$ws.Range("D1").Value = "automotiveRevenue"
$ws.Range("E1").Value = "automotiveCostOfRevenue"
$ws.Range("F1").Value = "deliveredModel3Y"
$ws.Range("G1").Value = "deliveredOtherModels"
$ws.Range("H1").Value = "producedModel3Y"
$ws.Range("I1").Value = "producedOtherModels"
$ws.Range("J1").Value = "producedCars"
$ws.Range("K1").Value = "deliveredCars"

# Fill in the new "producedCars" column with the production total
# (producedModel3Y + producedOtherModels), matching each row's style.
for ($r = 2; $r -le 21; $r++) {
    $ws.Range("J$r").Formula = "=H$r+I$r"
}
